$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") rows 2-29, replacing the old "Strike#"-derived values
$kValues = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 1
    6  = 1
    7  = 0
    8  = 0
    9  = 1
    10 = 2
    11 = 1
    12 = 1
    13 = 2
    14 = 0
    15 = 1
    16 = 2
    17 = 0
    18 = 1
    19 = 1
    20 = 3
    21 = 1
    22 = 2
    23 = 2
    24 = 1
    25 = 0
    26 = 1
    27 = 3
    28 = 1
    29 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
